$wb = $excel.ActiveWorkbook

# --- Sheet 1 (SHEET1) ---
$ws = $wb.Worksheets.Item(1)
$templateRow = 495

$r = 496
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "1-9 years"
$ws.Cells.Item($r, 3).Value = 461
$ws.Cells.Item($r, 4).Value = 1.8620243962
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 0

$r = 497
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "10-19 years"
$ws.Cells.Item($r, 3).Value = 1198
$ws.Cells.Item($r, 4).Value = 4.8388399709
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 0

$r = 498
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "20-29 years"
$ws.Cells.Item($r, 3).Value = 4538
$ws.Cells.Item($r, 4).Value = 18.329428871
$ws.Cells.Item($r, 5).Value = 4
$ws.Cells.Item($r, 6).Value = 0.5633802817

$r = 499
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "30-39 years"
$ws.Cells.Item($r, 3).Value = 4577
$ws.Cells.Item($r, 4).Value = 18.486953712
$ws.Cells.Item($r, 5).Value = 12
$ws.Cells.Item($r, 6).Value = 1.6901408451

$r = 500
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "40-49 years"
$ws.Cells.Item($r, 3).Value = 4407
$ws.Cells.Item($r, 4).Value = 17.800306971
$ws.Cells.Item($r, 5).Value = 26
$ws.Cells.Item($r, 6).Value = 3.661971831

$r = 501
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "50-59 years"
$ws.Cells.Item($r, 3).Value = 4146
$ws.Cells.Item($r, 4).Value = 16.74610227
$ws.Cells.Item($r, 5).Value = 75
$ws.Cells.Item($r, 6).Value = 10.563380282

$r = 502
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "60-64 years"
$ws.Cells.Item($r, 3).Value = 1650
$ws.Cells.Item($r, 4).Value = 6.6645124808
$ws.Cells.Item($r, 5).Value = 71
$ws.Cells.Item($r, 6).Value = 10

$r = 503
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "65-69 years"
$ws.Cells.Item($r, 3).Value = 1196
$ws.Cells.Item($r, 4).Value = 4.830761774
$ws.Cells.Item($r, 5).Value = 82
$ws.Cells.Item($r, 6).Value = 11.549295775

$r = 504
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "70-74 years"
$ws.Cells.Item($r, 3).Value = 778
$ws.Cells.Item($r, 4).Value = 3.1424186122
$ws.Cells.Item($r, 5).Value = 78
$ws.Cells.Item($r, 6).Value = 10.985915493

$r = 505
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "75-79 years"
$ws.Cells.Item($r, 3).Value = 586
$ws.Cells.Item($r, 4).Value = 2.3669117053
$ws.Cells.Item($r, 5).Value = 80
$ws.Cells.Item($r, 6).Value = 11.267605634

$r = 506
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "80+ years"
$ws.Cells.Item($r, 3).Value = 1103
$ws.Cells.Item($r, 4).Value = 4.455125616
$ws.Cells.Item($r, 5).Value = 282
$ws.Cells.Item($r, 6).Value = 39.718309859

$r = 507
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "<1 year"
$ws.Cells.Item($r, 3).Value = 102
$ws.Cells.Item($r, 4).Value = 0.4119880443
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 0

$r = 508
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "Unknown"
$ws.Cells.Item($r, 3).Value = 16
$ws.Cells.Item($r, 4).Value = 0.0646255756
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 0


# --- Sheet 2 (SHEET2) ---
$ws = $wb.Worksheets.Item(2)
$templateRow = 115

$r = 116
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "Female"
$ws.Cells.Item($r, 3).Value = 11928
$ws.Cells.Item($r, 4).Value = 48.178366589
$ws.Cells.Item($r, 5).Value = 292
$ws.Cells.Item($r, 6).Value = 41.126760563

$r = 117
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "Male"
$ws.Cells.Item($r, 3).Value = 12560
$ws.Cells.Item($r, 4).Value = 50.731076824
$ws.Cells.Item($r, 5).Value = 415
$ws.Cells.Item($r, 6).Value = 58.450704225

$r = 118
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "Unknown"
$ws.Cells.Item($r, 3).Value = 270
$ws.Cells.Item($r, 4).Value = 1.0905565878
$ws.Cells.Item($r, 5).Value = 3
$ws.Cells.Item($r, 6).Value = 0.4225352113


# --- Sheet 3 (SHEET3) ---
$ws = $wb.Worksheets.Item(3)
$templateRow = 223

$r = 224
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "Asian"
$ws.Cells.Item($r, 3).Value = 640
$ws.Cells.Item($r, 4).Value = 2.5850230229
$ws.Cells.Item($r, 5).Value = 14
$ws.Cells.Item($r, 6).Value = 1.9718309859

$r = 225
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "Black"
$ws.Cells.Item($r, 3).Value = 2851
$ws.Cells.Item($r, 4).Value = 11.515469747
$ws.Cells.Item($r, 5).Value = 91
$ws.Cells.Item($r, 6).Value = 12.816901408

$r = 226
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "Hispanic"
$ws.Cells.Item($r, 3).Value = 9712
$ws.Cells.Item($r, 4).Value = 39.227724372
$ws.Cells.Item($r, 5).Value = 205
$ws.Cells.Item($r, 6).Value = 28.873239437

$r = 227
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "Other"
$ws.Cells.Item($r, 3).Value = 132
$ws.Cells.Item($r, 4).Value = 0.5331609985
$ws.Cells.Item($r, 5).Value = 1
$ws.Cells.Item($r, 6).Value = 0.1408450704

$r = 228
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "Unknown"
$ws.Cells.Item($r, 3).Value = 5234
$ws.Cells.Item($r, 4).Value = 21.140641409
$ws.Cells.Item($r, 5).Value = 90
$ws.Cells.Item($r, 6).Value = 12.676056338

$r = 229
$ws.Range("A" + $templateRow + ":F" + $templateRow).Copy($ws.Range("A" + $r + ":F" + $r))
$ws.Cells.Item($r, 1).Value = 44022
$ws.Cells.Item($r, 2).Value = "White"
$ws.Cells.Item($r, 3).Value = 6189
$ws.Cells.Item($r, 4).Value = 24.997980451
$ws.Cells.Item($r, 5).Value = 309
$ws.Cells.Item($r, 6).Value = 43.521126761

